# Ground human oversight moderator in automation theory (Parasuraman et al., 2000)
# - Add an [INSTRUCTIONS] row (row 2) to the Study_Metadata sheet explaining the
#   oversight taxonomy collapse (10-level automation model -> 3 levels).
# - Add a cell comment on C1 documenting the Human Oversight Level taxonomy,
#   grounded in Parasuraman et al. (2000) and the EU AI Act (2024, Art. 14).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Study_Metadata")

# --- New instructions row under the header row ---------------------------
$ws.Range("A2").Value = "[INSTRUCTIONS]"
$ws.Range("Q2").Value = "Oversight taxonomy: Parasuraman et al. (2000) 10-level automation model collapsed to 3 levels. See AI_Agent_Characteristics sheet, column C."

# --- Cell comment on C1 grounding the oversight taxonomy ------------------
$commentText = @"
Human Oversight Level (Parasuraman et al., 2000)
Three-level taxonomy grounded in the 10-level automation model:
  fully_autonomous = Automation levels 7-10 (computer decides, may/may not inform human)
  ai_led_checkpoints = Automation levels 4-6 (computer suggests/executes, human can intervene)
  human_led_ai_support = Automation levels 1-3 (human decides with computer assistance)
See coding_manual.md §C1 for full decision rules.
Policy context: EU AI Act (2024, Art. 14) mandates human oversight for high-risk AI in education.
"@

$comment = $ws.Range("C1").AddComment($commentText)
$comment.Author = "Meta-Analysis Coding Team"
$comment.Visible = $false
